$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header columns: CLASS NAME -> MODULE, PACKAGE SUFFIX -> SUBMODULE
$ws.Range("C1").Value = "MODULE"
$ws.Range("D1").Value = "SUBMODULE"

# Replace CLASS NAME / PACKAGE SUFFIX data with MODULE / SUBMODULE values
$ws.Range("C2").Value = "Nursery"
$ws.Range("D2").Value = "Germination"

$ws.Range("C3").Value = "Nursery"
$ws.Range("D3").Value = "Germination"

$ws.Range("C4").Value = "Nursery"
$ws.Range("D4").Value = "Setup"

# Drop the now-unused IFS PROJECTION column (E)
$ws.Range("E1:E4").Clear()

# Move the active selection like the recorded session
$ws.Range("D16").Select()
